# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holdings detail) right before the
#    "总计" (totals) summary sheet.
# 2. Insert a new leading row into "总计" for the 2022-Q1 quarter and
#    renumber the existing index column.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New "2022-Q1" sheet, placed after "2021-Q4" (i.e. right before "总计").
# ---------------------------------------------------------------------------
$afterSheet = $wb.Worksheets.Item("2021-Q4")
$q1 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $afterSheet)
$q1.Name = "2022-Q1"

$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $q1.Cells.Item(1, 2 + $i)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Columns: 基金代码, 基金名称, 基金规模, 股票总仓位, 仓位占比, 持有市值(亿元) (text) then 仓位排名 (number)
$rows = @(
    @("970041", "国海量化优选一年持有股票A", "7.70", "91.93", "0.72", "0.0554", 7),
    @("970042", "国海量化优选一年持有股票C", "6.95", "91.93", "0.72", "0.0500", 7),
    @("009847", "圆信永丰研究精选混合A",     "1.83", "90.78", "2.62", "0.0479", 10),
    @("006969", "圆信永丰高端制造混合",       "1.04", "86.60", "3.44", "0.0358", 7),
    @("290008", "泰信发展主题混合",           "0.68", "81.03", "5.13", "0.0349", 5),
    @("000270", "建信灵活配置混合",           "2.12", "93.93", "0.87", "0.0184", 10),
    @("009848", "圆信永丰研究精选混合C",     "0.46", "90.78", "2.62", "0.0121", 10),
    @("005247", "国都量化精选混合",           "0.02", "64.74", "2.13", "0.0004", 8)
)

$r = 2
foreach ($row in $rows) {
    $idxCell = $q1.Cells.Item($r, 1)
    $idxCell.Value = $r - 2
    $idxCell.Font.Bold = $true
    $idxCell.Borders.LineStyle = 1
    $idxCell.HorizontalAlignment = -4108
    $idxCell.VerticalAlignment = -4160

    for ($c = 0; $c -lt 6; $c++) {
        $textCell = $q1.Cells.Item($r, 2 + $c)
        $textCell.NumberFormat = "@"
        $textCell.Value = $row[$c]
    }
    $q1.Cells.Item($r, 8).Value = $row[6]

    $r = $r + 1
}

# ---------------------------------------------------------------------------
# 2) Update "总计": insert the new 2022-Q1 row at the top of the data and
#    renumber the existing rows' index column.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()

$newIdx = $total.Cells.Item(2, 1)
$newIdx.Value = 0
$newIdx.Font.Bold = $true
$newIdx.Borders.LineStyle = 1
$newIdx.HorizontalAlignment = -4108
$newIdx.VerticalAlignment = -4160

$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 8
$total.Cells.Item(2, 4).Value = 0.25

# Renumber the index column for the rows that got pushed down (they kept
# their old 0,1,2 values; they should now read 1,2,3).
for ($row = 3; $row -le 5; $row++) {
    $total.Cells.Item($row, 1).Value = $row - 2
}
